# Apply "repull data, push all data, mean calculation" edit:
# update the dSF column (column F) values for several rows to reflect
# freshly pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 4
$ws.Range("F4").Value = -1
$ws.Range("F5").Value = 2
$ws.Range("F9").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("F16").Value = -3
